# Apply updated cryptocurrency price/volume data to Sheet1.
# Values in column D that could be parsed as plain numbers are prefixed
# with a leading apostrophe so Excel stores them as text (preserving
# formatting such as trailing zeros, e.g. "0.05330"), matching the
# original workbook's inline-string / text-typed cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.913.87"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "1.742.33"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'231.19"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").Value = "'0.9992"
$ws.Range("D7").Value = "'0.5266"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "'0.2773"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "'39.57"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").Value = "'0.06157"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "1.734.62"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").Value = "'0.07125"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "'15.26"
$ws.Range("E13").Value = "  -3.30%  "
$ws.Range("D14").Value = "'0.6471"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "'77.22"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "'0.9994"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "25.878.13"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").Value = "'11.56"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'0.000006676"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "1.959.05"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").Value = "'4.266"
$ws.Range("E23").Value = "  +4.73%  "
$ws.Range("D24").Value = "'8.806"
$ws.Range("E24").Value = "  +4.13%  "
$ws.Range("D25").Value = "'5.181"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'140.57"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("D28").Value = "'15.24"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "'1.806"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("D30").Value = "'102.51"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "'0.08339"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").Value = "'3.741"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").Value = "'3.595"
$ws.Range("E33").Value = "  +4.10%  "
$ws.Range("D34").Value = "'0.04516"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "'2.610"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "'0.9797"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "'0.6226"
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("D38").Value = "'2.688"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "'0.01586"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "'1.931"
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").Value = "'0.9989"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "'100.16"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("D43").Value = "'0.3878"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'0.7314"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "'5.026"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("D46").Value = "'0.05330"
$ws.Range("E46").Value = "  -3.40%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "'6.253"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").Value = "'53.68"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'30.17"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'7.656"
$ws.Range("E51").Value = "  +2.97%  "
